$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (Fri Nov 15 17:53:54 UTC 2024 refresh)
$ws.Range("D2").Value = "89.614.67"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "3.045.74"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.24"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "611.96"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.360"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.868"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.12%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "3.045.68"
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.675"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +21.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.187"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.37"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "89.361.59"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.33"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "3.046.40"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000220"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.17"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.24"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.02"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.39"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.85"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.61"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.35"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.74"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "502.24"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.65"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.77"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.56%  "
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("E38").Value = "  -10.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.23"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.363"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.83"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.55"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.42"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0687"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.19"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "161.61"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.697"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.29%  "
